$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Move "Straight Connector 45" (id 46) ---
$conn45 = $s.Shapes.Item(33)
$conn45.Left = 502.6319
$conn45.Top  = 115.8993

# --- Move "Straight Arrow Connector 47" (id 48) ---
$conn47 = $s.Shapes.Item(36)
$conn47.Left = 505.8196
$conn47.Top  = 199.89433

# --- Move "Straight Arrow Connector 144" (id 145) ---
$conn144 = $s.Shapes.Item(64)
$conn144.Left = 507.9698
$conn144.Top  = 240.19182

# --- Add new "Straight Arrow Connector 73" dashed purple arrow ---
# Duplicate a sibling connector so the <p:style> block / cxnSpLocks match
# the authored shape exactly, then re-position & re-style it.
$newConnShape = $conn144.Duplicate()
$newConn = $s.Shapes.Item($s.Shapes.Count)
$newConn.Name = "Straight Arrow Connector 73"

$newConn.Left   = 507.9698
$newConn.Top    = 253.8901
$newConn.Width  = 100.9941
$newConn.Height = 0.3523
$newConn.Flip(1)

$newConn.Line.Weight = 1.5
$newConn.Line.ForeColor.RGB = 0xA03070
$newConn.Line.DashStyle = 9
$newConn.Line.BeginArrowheadStyle = 3
$newConn.Line.BeginArrowheadWidth = 2
$newConn.Line.BeginArrowheadLength = 2
$newConn.Line.EndArrowheadStyle = 1
$newConn.Line.EndArrowheadWidth = 2
$newConn.Line.EndArrowheadLength = 2
